# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" timestamps for the
#   861d4a25-7293-4eca-8e95-1c50880602b1.md item (Overview + de-de sheets
#   share one timestamp value, zh-cn sheet has its own).
# - Sets the "Priority" column to "ht" for the six items that previously
#   had no priority set, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 10, 11, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-19 20:17:53"
}

# zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-08-19 20:17:47"
    $wsZh.Range("E$r").Value = "ht"
}

# de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority"
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-08-19 20:17:53"
    $wsDe.Range("E$r").Value = "ht"
}
